$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column O ("Unnamed: 14") is an empty placeholder column with no data rows.
# Delete it so column P ("Imágenes" header + the image-link values) shifts
# left into column O.
$ws.Columns("O").Delete()
